# Resume edit:
#  1. Drop the stray _GoBack bookmark that sat at the very top of the
#     document (wrapping nothing, right before the "Education" heading).
#  2. The author's cursor (_GoBack) ended up instead in the middle of the
#     "... experienced problems using Moolah's services" bullet, right
#     after "Moolah " - and "'s" was deleted so it now reads
#     "... Moolah services". That forces the run that used to read
#     "oblems using Moolah's services" to split in two, with the
#     (moved) _GoBack bookmark sitting between them.

$d = $word.ActiveDocument

# --- Step 1: remove the old _GoBack bookmark near the top of the doc ---
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Step 2: find the bullet paragraph that needs the text/bookmark edit ---
$needle = "oblems using Moolah" + [char]0x2019 + "s services"
$hit = $d.Content
$hit.Find.Execute($needle) | Out-Null

# Build the corrected paragraph: identical pPr/formatting to the original,
# but with "'s" removed and the run split exactly where the cursor (the
# _GoBack bookmark) now sits, between "Moolah " and "services".
$newParagraphXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6B5B52F1" w14:textId="77777777" w:rsidR="00073751" w:rsidRPr="00823FC0" w:rsidRDefault="00073751" w:rsidP="00073751"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="20" w:lineRule="atLeast"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Applied previous interpersonal and customer service skills to provide</w:t></w:r><w:r w:rsidRPr="00823FC0"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> customer support to clients who experienced pr</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">oblems using Moolah </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>services</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Insert the fixed-up paragraph right after the old (still intact) one.
# Using InsertXML (rather than editing the run's Text in place) keeps the
# run split crisp instead of letting same-formatted neighbour runs coalesce.
$insertionPoint = $d.Range($hit.End, $hit.End)
$insertionPoint.InsertXML($newParagraphXml) | Out-Null

# --- Step 3: delete the now-duplicated original paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*" + $needle + "*") {
        $para.Range.Delete()
        break
    }
}
